$ws = $excel.ActiveWorkbook.ActiveSheet

# 1) Copy formatting (A:M) from the last existing row (470) into the new rows (471:490).
#    Column N is handled separately below since row 470 has no N cell to copy from.
$ws.Range("A470:M470").Copy()
$ws.Range("A471:M490").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("A471:A490").RowHeight = 15.75

# 2) Apply the same text style to column N only for the rows that actually use it,
#    so untouched rows do not get a stray empty N cell (mirrors source data: [473, 474, 475, 476, 478, 483, 485, 486, 489, 490]).
$ws.Range("M470").Copy()
$nRows = @(473, 474, 475, 476, 478, 483, 485, 486, 489, 490)
foreach ($r in $nRows) {
    $ws.Range("N$r").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

# 3) Fill in the actual survey-response values for the 20 new rows (471:490).

# Row 471
$ws.Range("A471").Value = 45193.82327327546
$ws.Range("B471").Value = "dlruddk9@naver.com"
$ws.Range("C471").Value = "사회복지학과"
$ws.Range("D471").Value = 20212342
$ws.Range("E471").Value = "이경아"
$ws.Range("F471").Value = "‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다."
$ws.Range("G471").Value = 0.5
$ws.Range("H471").Value = "5:5"
$ws.Range("I471").Value = "20분의 1"
$ws.Range("J471").Value = "44만호, 153만명"
$ws.Range("K471").Value = "전라"
$ws.Range("L471").Value = "Red"
$ws.Range("M471").Value = "반대한다."

# Row 472
$ws.Range("A472").Value = 45193.824795555556
$ws.Range("B472").Value = "ayden0429@gmail.com"
$ws.Range("C472").Value = "의예과"
$ws.Range("D472").Value = 20226145
$ws.Range("E472").Value = "이성연"
$ws.Range("F472").Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Range("G472").Value = 0.1
$ws.Range("H472").Value = "6:4"
$ws.Range("I472").Value = "20분의 1"
$ws.Range("J472").Value = "20만호, 69만명"
$ws.Range("K472").Value = "충청"
$ws.Range("L472").Value = "Red"
$ws.Range("M472").Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

# Row 473
$ws.Range("A473").Value = 45193.830308414355
$ws.Range("B473").Value = "lcbat4@gmail.com"
$ws.Range("C473").Value = "글로벌학부"
$ws.Range("D473").Value = 20236429
$ws.Range("E473").Value = "홍서경"
$ws.Range("F473").Value = "과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다."
$ws.Range("G473").Value = 0.7
$ws.Range("H473").Value = "6:4"
$ws.Range("I473").Value = "20분의 1"
$ws.Range("J473").Value = "20만호, 69만명"
$ws.Range("K473").Value = "경상"
$ws.Range("L473").Value = "Black"
$ws.Range("N473").Value = "모름/무응답"

# Row 474
$ws.Range("A474").Value = 45193.832151145834
$ws.Range("B474").Value = "hkmcosmos1@gmail.com"
$ws.Range("C474").Value = "글로벌비즈니스"
$ws.Range("D474").Value = 20226429
$ws.Range("E474").Value = "한기민"
$ws.Range("F474").Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Range("G474").Value = 0.1
$ws.Range("H474").Value = "3:7"
$ws.Range("I474").Value = "15분의 1"
$ws.Range("J474").Value = "20만호, 69만명"
$ws.Range("K474").Value = "충청"
$ws.Range("L474").Value = "Black"
$ws.Range("N474").Value = "모름/무응답"

# Row 475
$ws.Range("A475").Value = 45193.837822326386
$ws.Range("B475").Value = "kimbitna7890@naver.com"
$ws.Range("C475").Value = "광고홍보학과"
$ws.Range("D475").Value = 20222609
$ws.Range("E475").Value = "김빛나"
$ws.Range("F475").Value = "과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다."
$ws.Range("G475").Value = 0.3
$ws.Range("H475").Value = "3:7"
$ws.Range("I475").Value = "15분의 1"
$ws.Range("J475").Value = "15만호,  32만명"
$ws.Range("K475").Value = "평안"
$ws.Range("L475").Value = "Black"
$ws.Range("N475").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

# Row 476
$ws.Range("A476").Value = 45193.83826665509
$ws.Range("B476").Value = "hyelinj27@gmail.com"
$ws.Range("C476").Value = "인공지능융합학부"
$ws.Range("D476").Value = 20236781
$ws.Range("E476").Value = "진혜린"
$ws.Range("F476").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G476").Value = 0.1
$ws.Range("H476").Value = "6:4"
$ws.Range("I476").Value = "20분의 1"
$ws.Range("J476").Value = "20만호, 69만명"
$ws.Range("K476").Value = "충청"
$ws.Range("L476").Value = "Black"
$ws.Range("N476").Value = "찬성한다."

# Row 477
$ws.Range("A477").Value = 45193.84510828704
$ws.Range("B477").Value = "sillysunny@naver.com"
$ws.Range("C477").Value = "인문학부"
$ws.Range("D477").Value = 20231037
$ws.Range("E477").Value = "박세현"
$ws.Range("F477").Value = "‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다."
$ws.Range("G477").Value = 0.7
$ws.Range("H477").Value = "6:4"
$ws.Range("I477").Value = "20분의 1"
$ws.Range("J477").Value = "20만호, 69만명"
$ws.Range("K477").Value = "전라"
$ws.Range("L477").Value = "Red"
$ws.Range("M477").Value = "모름/무응답"

# Row 478
$ws.Range("A478").Value = 45193.852299907405
$ws.Range("B478").Value = "ljh2017@naver.com"
$ws.Range("C478").Value = "철학과"
$ws.Range("D478").Value = 20181079
$ws.Range("E478").Value = "이정효"
$ws.Range("F478").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G478").Value = 0.3
$ws.Range("H478").Value = "3:7"
$ws.Range("I478").Value = "10분의 1"
$ws.Range("J478").Value = "44만호, 153만명"
$ws.Range("K478").Value = "전라"
$ws.Range("L478").Value = "Black"
$ws.Range("N478").Value = "모름/무응답"

# Row 479
$ws.Range("A479").Value = 45193.85701881944
$ws.Range("B479").Value = "jisung5549@naver.com"
$ws.Range("C479").Value = "경영학과"
$ws.Range("D479").Value = 20222970
$ws.Range("E479").Value = "송지성"
$ws.Range("F479").Value = "실제로 현장에 나가서 수확량을 파악하고 등급을 매기는 답험(踏驗)을 하였다."
$ws.Range("G479").Value = 0.7
$ws.Range("H479").Value = "4:6"
$ws.Range("I479").Value = "20분의 1"
$ws.Range("J479").Value = "44만호, 153만명"
$ws.Range("K479").Value = "경상"
$ws.Range("L479").Value = "Red"
$ws.Range("M479").Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

# Row 480
$ws.Range("A480").Value = 45193.862970636575
$ws.Range("B480").Value = "jangho5636@gmail.com"
$ws.Range("C480").Value = "러시아학과"
$ws.Range("D480").Value = 20161723
$ws.Range("E480").Value = "이장호"
$ws.Range("F480").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G480").Value = 0.9
$ws.Range("H480").Value = "7:3"
$ws.Range("I480").Value = "30분의 1"
$ws.Range("J480").Value = "20만호, 69만명"
$ws.Range("K480").Value = "평안"
$ws.Range("L480").Value = "Red"
$ws.Range("M480").Value = "반대한다."

# Row 481
$ws.Range("A481").Value = 45193.8716166088
$ws.Range("B481").Value = "0227jsh@naver.com"
$ws.Range("C481").Value = "식품영양학과"
$ws.Range("D481").Value = 20233843
$ws.Range("E481").Value = "장서희"
$ws.Range("F481").Value = "과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다."
$ws.Range("G481").Value = 0.3
$ws.Range("H481").Value = "6:4"
$ws.Range("I481").Value = "20분의 1"
$ws.Range("J481").Value = "15만호,  32만명"
$ws.Range("K481").Value = "평안"
$ws.Range("L481").Value = "Red"
$ws.Range("M481").Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

# Row 482
$ws.Range("A482").Value = 45193.8724659375
$ws.Range("B482").Value = "rhksan324@naver.com"
$ws.Range("C482").Value = "금융재무학과"
$ws.Range("D482").Value = 20203001
$ws.Range("E482").Value = "이관무"
$ws.Range("F482").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G482").Value = 0.1
$ws.Range("H482").Value = "6:4"
$ws.Range("I482").Value = "20분의 1"
$ws.Range("J482").Value = "20만호, 69만명"
$ws.Range("K482").Value = "충청"
$ws.Range("L482").Value = "Red"
$ws.Range("M482").Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

# Row 483
$ws.Range("A483").Value = 45193.8728371412
$ws.Range("B483").Value = "jihye199530@gmail.com"
$ws.Range("C483").Value = "간호학과"
$ws.Range("D483").Value = 20236261
$ws.Range("E483").Value = "엄지혜"
$ws.Range("F483").Value = "‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다."
$ws.Range("G483").Value = 0.7
$ws.Range("H483").Value = "5:5"
$ws.Range("I483").Value = "20분의 1"
$ws.Range("J483").Value = "15만호,  32만명"
$ws.Range("K483").Value = "평안"
$ws.Range("L483").Value = "Black"
$ws.Range("N483").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

# Row 484
$ws.Range("A484").Value = 45193.87713403935
$ws.Range("B484").Value = "jongbaep17s@gmail.com"
$ws.Range("C484").Value = "글로벌비즈니스"
$ws.Range("D484").Value = 20226410
$ws.Range("E484").Value = "박종배"
$ws.Range("F484").Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Range("G484").Value = 0.1
$ws.Range("H484").Value = "6:4"
$ws.Range("I484").Value = "20분의 1"
$ws.Range("J484").Value = "20만호, 69만명"
$ws.Range("K484").Value = "충청"
$ws.Range("L484").Value = "Red"
$ws.Range("M484").Value = "모름/무응답"

# Row 485
$ws.Range("A485").Value = 45193.9036634838
$ws.Range("B485").Value = "rhdskrud123@naver.com"
$ws.Range("C485").Value = "인문학부"
$ws.Range("D485").Value = 20231002
$ws.Range("E485").Value = "공나경"
$ws.Range("F485").Value = "실제로 현장에 나가서 수확량을 파악하고 등급을 매기는 답험(踏驗)을 하였다."
$ws.Range("G485").Value = 0.7
$ws.Range("H485").Value = "7:3"
$ws.Range("I485").Value = "20분의 1"
$ws.Range("J485").Value = "20만호, 69만명"
$ws.Range("K485").Value = "전라"
$ws.Range("L485").Value = "Black"
$ws.Range("N485").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

# Row 486
$ws.Range("A486").Value = 45193.90699703703
$ws.Range("B486").Value = "shdbsgh0305@naver.com"
$ws.Range("C486").Value = "러시아학과"
$ws.Range("D486").Value = 20231710
$ws.Range("E486").Value = "노윤호"
$ws.Range("F486").Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Range("G486").Value = 0.5
$ws.Range("H486").Value = "5:5"
$ws.Range("I486").Value = "20분의 1"
$ws.Range("J486").Value = "20만호, 69만명"
$ws.Range("K486").Value = "전라"
$ws.Range("L486").Value = "Black"
$ws.Range("N486").Value = "모름/무응답"

# Row 487
$ws.Range("A487").Value = 45193.90846236111
$ws.Range("B487").Value = "digiphk12@naver.com"
$ws.Range("C487").Value = "광고홍보학과"
$ws.Range("D487").Value = 20162617
$ws.Range("E487").Value = "박현규"
$ws.Range("F487").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G487").Value = 0.7
$ws.Range("H487").Value = "6:4"
$ws.Range("I487").Value = "15분의 1"
$ws.Range("J487").Value = "20만호, 69만명"
$ws.Range("K487").Value = "평안"
$ws.Range("L487").Value = "Red"
$ws.Range("M487").Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

# Row 488
$ws.Range("A488").Value = 45193.909824409726
$ws.Range("B488").Value = "mnsghn314@naver.com"
$ws.Range("C488").Value = "소프트웨어학과"
$ws.Range("D488").Value = 20235159
$ws.Range("E488").Value = "문승현"
$ws.Range("F488").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G488").Value = 0.1
$ws.Range("H488").Value = "6:4"
$ws.Range("I488").Value = "10분의 1"
$ws.Range("J488").Value = "20만호, 69만명"
$ws.Range("K488").Value = "충청"
$ws.Range("L488").Value = "Red"
$ws.Range("M488").Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

# Row 489
$ws.Range("A489").Value = 45193.91182700232
$ws.Range("B489").Value = "amielee1997@naver.com"
$ws.Range("C489").Value = "간호학과"
$ws.Range("D489").Value = 20217159
$ws.Range("E489").Value = "이지수"
$ws.Range("F489").Value = "과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다."
$ws.Range("G489").Value = 0.1
$ws.Range("H489").Value = "7:3"
$ws.Range("I489").Value = "10분의 1"
$ws.Range("J489").Value = "15만호,  32만명"
$ws.Range("K489").Value = "경기"
$ws.Range("L489").Value = "Black"
$ws.Range("N489").Value = "모름/무응답"

# Row 490
$ws.Range("A490").Value = 45193.91356240741
$ws.Range("B490").Value = "sinfkks@gmail.com"
$ws.Range("C490").Value = "반도체 디스플레이스쿨"
$ws.Range("D490").Value = 20233304
$ws.Range("E490").Value = "김경진"
$ws.Range("F490").Value = "‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다."
$ws.Range("G490").Value = 0.9
$ws.Range("H490").Value = "4:6"
$ws.Range("I490").Value = "10분의 1"
$ws.Range("J490").Value = "44만호, 153만명"
$ws.Range("K490").Value = "경기"
$ws.Range("L490").Value = "Black"
$ws.Range("N490").Value = "찬성한다."

# 4) Match the saved view state (active cell after entering the new rows).
$ws.Range("G499").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 470
